# Append five new "Paragraphedeliste"-styled paragraphs at the end of the
# document (after the existing "Lucas ... seq " bullet item and before the
# section properties), as per the carnet de bord update:
#
#   - "Powerpoint fini 3h"
#   - (empty)
#   - "03/05/2021"
#   - "On à terminer les choix de la solution. On commence à taper dans le dur"
#   - (empty)

$d = $word.ActiveDocument

function Add-CarnetParagraph($Text) {
    # Grow the document by one paragraph at the very end of the story.
    $tail = $d.Paragraphs.Last.Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()

    $para = $d.Paragraphs.Last
    # Force the "Paragraphedeliste" style explicitly on the new paragraph so
    # it does not merely inherit the previous paragraph's numbered-list
    # formatting (numPr) -- the new entries are plain list-style paragraphs,
    # not numbered items.
    $para.Style = "Paragraphedeliste"

    if ($Text -ne "") {
        $para.Range.Text = $Text
    } else {
        # Leave a genuinely empty paragraph (no run/text) -- insert then
        # remove a placeholder character so no stray run is left behind.
        $para.Range.InsertAfter("x")
        $para.Range.Characters(1).Delete()
    }

    return $para
}

Add-CarnetParagraph "Powerpoint fini 3h" | Out-Null
Add-CarnetParagraph "" | Out-Null
Add-CarnetParagraph "03/05/2021" | Out-Null
Add-CarnetParagraph "On à terminer les choix de la solution. On commence à taper dans le dur" | Out-Null
Add-CarnetParagraph "" | Out-Null

Write-Output "Added 5 paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
